# Fix death on treatment MDR default data
# Death rate too high for MDR deaths on treatment.
#
# The "program_perc_treatment_death_mdr" row (row 2) on the time_variants
# sheet is removed entirely; the rows below it (program_perc_detect,
# int_perc_ipt_age0to5) shift up to take its place.

$wb = $excel.ActiveWorkbook

$wsTime = $wb.Worksheets.Item("time_variants")
$wsConst = $wb.Worksheets.Item("constants")

# Remember the smoothness-column validation so it can be restored with its
# original full-column extent after the row shift below.
$smoothRange = $wsTime.Range("C2:C1048576")

# Delete the whole row holding program_perc_treatment_death_mdr; this shifts
# program_perc_detect and int_perc_ipt_age0to5 up one row each.
$wsTime.Rows.Item(2).Delete()

# Restore the smoothness-column validation to its original extent (row
# deletion shrinks a range that ends at the sheet's last row by one row).
$smoothRange.Validation.Delete()
$smoothRange.Validation.Add(2, 1, 1, "0", "100")
$smoothRange.Validation.InputTitle = "Smoothness for fitting function"
$smoothRange.Validation.InputMessage = "Must be positive."
$smoothRange.Validation.ShowInput = $true
$smoothRange.Validation.ShowError = $true

# Selection on constants sheet moves to entire row 5.
$wsConst.Rows.Item(5).Select()

# Selection (bottom-right frozen pane) on time_variants moves to A9,
# and that sheet becomes the active tab.
$wsTime.Range("A9").Select()
$wsTime.Activate()
